# Scheduled-runner refresh of FFXIV leve-crafting profit figures (Exodus
# server) across all eight Disciple of the Hand sheets. Updated market
# price / profit columns (H:N) for the specific leves whose prices moved;
# a couple of rows also gain/lose a profit cell (NQ vs HQ) that previously
# didn't apply.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3361.4
$ws.Range("I86").Value = 2947.6667
$ws.Range("J86").Value = 3982
$ws.Range("K86").Value = 2947.6667
$ws.Range("L86").Value = 3982
$ws.Range("M86").Value = -1824.6667
$ws.Range("N86").Value = -6228

$ws.Range("H89").Value = 3361.4
$ws.Range("I89").Value = 2947.6667
$ws.Range("J89").Value = 3982
$ws.Range("K89").Value = 14738.3335
$ws.Range("L89").Value = 19910
$ws.Range("M89").Value = -9122.333500000001
$ws.Range("N89").Value = -31142

$ws.Range("H113").Value = 4883.1665
$ws.Range("I113").Value = 4774.75
$ws.Range("J113").Value = 5100
$ws.Range("K113").Value = 4774.75
$ws.Range("L113").Value = 5100
$ws.Range("M113").Value = -1520.75
$ws.Range("N113").Value = -11608

$ws.Range("H125").Value = 506013
$ws.Range("I125").Value = 506013
$ws.Range("K125").Value = 4554117
$ws.Range("M125").Value = -4551657

$ws.Range("H141").Value = 2544
$ws.Range("I141").Value = 2214.3333
$ws.Range("K141").Value = 6642.999899999999
$ws.Range("M141").Value = -1462.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 858.1429000000001
$ws.Range("I2").Value = 873.5
$ws.Range("J2").Value = 837.6667
$ws.Range("K2").Value = 873.5
$ws.Range("L2").Value = 837.6667
$ws.Range("M2").Value = -760.5
$ws.Range("N2").Value = -1063.6667

$ws.Range("H32").Value = 5767.8394
$ws.Range("I32").Value = 2203.6304
$ws.Range("K32").Value = 2203.6304
$ws.Range("M32").Value = -1916.6304

$ws.Range("H116").Value = 858.1429000000001
$ws.Range("I116").Value = 873.5
$ws.Range("J116").Value = 837.6667
$ws.Range("K116").Value = 873.5
$ws.Range("L116").Value = 837.6667
$ws.Range("M116").Value = 1420.5
$ws.Range("N116").Value = -5425.6667

$ws.Range("H122").Value = 4181.857
$ws.Range("I122").Value = 3344.2222
$ws.Range("J122").Value = 5689.6
$ws.Range("K122").Value = 10032.6666
$ws.Range("L122").Value = 17068.8
$ws.Range("M122").Value = -7582.6666
$ws.Range("N122").Value = -21968.8

$ws.Range("H132").Value = 1990.9
$ws.Range("I132").Value = 1990.9
$ws.Range("K132").Value = 5972.700000000001
$ws.Range("M132").Value = -3442.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 858.1429000000001
$ws.Range("I3").Value = 873.5
$ws.Range("J3").Value = 837.6667
$ws.Range("K3").Value = 873.5
$ws.Range("L3").Value = 837.6667
$ws.Range("M3").Value = -759.5
$ws.Range("N3").Value = -1065.6667

$ws.Range("H22").Value = 101354.5
$ws.Range("I22").Value = 167797.17
$ws.Range("K22").Value = 167797.17
$ws.Range("M22").Value = -167624.17

$ws.Range("H86").Value = 4574.8335
$ws.Range("I86").Value = 4011.8
$ws.Range("K86").Value = 4011.8
$ws.Range("M86").Value = -2888.8

$ws.Range("H89").Value = 4574.8335
$ws.Range("I89").Value = 4011.8
$ws.Range("K89").Value = 20059
$ws.Range("M89").Value = -14443

$ws.Range("H96").Value = 44444
$ws.Range("I96").Value = 44444
$ws.Range("K96").Value = 44444
$ws.Range("M96").Value = -41698

$ws.Range("H105").Value = 53341.45
$ws.Range("I105").Value = 69392.60000000001
$ws.Range("J105").Value = 5188
$ws.Range("K105").Value = 69392.60000000001
$ws.Range("L105").Value = 5188
$ws.Range("M105").Value = -67645.60000000001
$ws.Range("N105").Value = -8682

$ws.Range("H110").Value = 78922
$ws.Range("J110").Value = 78922
$ws.Range("L110").Value = 78922
$ws.Range("N110").Value = -87102

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2165
$ws.Range("I16").Value = 1998
$ws.Range("K16").Value = 1998
$ws.Range("M16").Value = -1711

$ws.Range("H58").Value = 1634.7273
$ws.Range("J58").Value = 2352.25
$ws.Range("L58").Value = 2352.25
$ws.Range("N58").Value = -2758.25

$ws.Range("H111").Value = 69999.5
$ws.Range("J111").Value = 69999.5
$ws.Range("L111").Value = 69999.5
$ws.Range("N111").Value = -78179.5

$ws.Range("H113").Value = 2165
$ws.Range("I113").Value = 1998
$ws.Range("K113").Value = 1998
$ws.Range("M113").Value = 172

$ws.Range("H122").Value = 3051.6667
$ws.Range("I122").Value = 2273.375
$ws.Range("J122").Value = 3674.3
$ws.Range("K122").Value = 6820.125
$ws.Range("L122").Value = 11022.9
$ws.Range("M122").Value = -4370.125
$ws.Range("N122").Value = -15922.9

$ws.Range("H134").Value = 2500091.5
$ws.Range("I134").Value = 2978808.2
$ws.Range("J134").Value = 202251.6
$ws.Range("K134").Value = 8936424.600000001
$ws.Range("L134").Value = 606754.8
$ws.Range("M134").Value = -8933889.600000001
$ws.Range("N134").Value = -611824.8

$ws.Range("H136").Value = 1634.7273
$ws.Range("J136").Value = 2352.25
$ws.Range("L136").Value = 7056.75
$ws.Range("N136").Value = -12156.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5299.6665
$ws.Range("I56").Value = 5299.6665
$ws.Range("K56").Value = 5299.6665
$ws.Range("M56").Value = -4769.6665

$ws.Range("H68").Value = 336333
$ws.Range("J68").Value = 503499.5
$ws.Range("L68").Value = 1510498.5
$ws.Range("N68").Value = -1512120.5

$ws.Range("H71").Value = 336333
$ws.Range("J71").Value = 503499.5
$ws.Range("L71").Value = 4531495.5
$ws.Range("N71").Value = -4539607.5

$ws.Range("H97").Value = 107.92308
$ws.Range("J97").Value = 104.4
$ws.Range("L97").Value = 313.2
$ws.Range("N97").Value = -1305.2

$ws.Range("H117").Value = 735.7778
$ws.Range("I117").Value = 657
$ws.Range("J117").Value = 798.8
$ws.Range("K117").Value = 1971
$ws.Range("L117").Value = 2396.4
$ws.Range("M117").Value = 1471
$ws.Range("N117").Value = -9280.4

$ws.Range("H121").Value = 2767.5557
$ws.Range("J121").Value = 3666.5
$ws.Range("L121").Value = 10999.5
$ws.Range("N121").Value = -13619.5

$ws.Range("H131").Value = 1433.5
$ws.Range("I131").Value = 862.5
$ws.Range("J131").Value = 2194.8333
$ws.Range("K131").Value = 2587.5
$ws.Range("L131").Value = 6584.499899999999
$ws.Range("M131").Value = 2452.5
$ws.Range("N131").Value = -16664.4999

$ws.Range("H132").Value = 3146.4814
$ws.Range("I132").Value = 683.2222
$ws.Range("J132").Value = 4378.1113
$ws.Range("K132").Value = 6148.999800000001
$ws.Range("L132").Value = 39403.00169999999
$ws.Range("M132").Value = -3618.999800000001
$ws.Range("N132").Value = -44463.00169999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9365.556
$ws.Range("I70").Value = 11330
$ws.Range("J70").Value = 8383.333000000001
$ws.Range("K70").Value = 11330
$ws.Range("L70").Value = 8383.333000000001
$ws.Range("M70").Value = -11060
$ws.Range("N70").Value = -8923.333000000001

$ws.Range("H73").Value = 9365.556
$ws.Range("I73").Value = 11330
$ws.Range("J73").Value = 8383.333000000001
$ws.Range("K73").Value = 11330
$ws.Range("L73").Value = 8383.333000000001
$ws.Range("M73").Value = -10394
$ws.Range("N73").Value = -10255.333

$ws.Range("H103").Value = 29680
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 29680
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 29680
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -32024

$ws.Range("H113").Value = 3231.6
$ws.Range("I113").Value = 2969.3333
$ws.Range("K113").Value = 2969.3333
$ws.Range("M113").Value = -799.3332999999998

$ws.Range("H122").Value = 6252996.5
$ws.Range("I122").Value = 10002745
$ws.Range("J122").Value = 3416.3333
$ws.Range("K122").Value = 30008235
$ws.Range("L122").Value = 10248.9999
$ws.Range("M122").Value = -30005785
$ws.Range("N122").Value = -15148.9999

$ws.Range("H126").Value = 3606.7778
$ws.Range("I126").Value = 2474.875
$ws.Range("J126").Value = 4512.3
$ws.Range("K126").Value = 7424.625
$ws.Range("L126").Value = 13536.9
$ws.Range("M126").Value = -4954.625
$ws.Range("N126").Value = -18476.9

$ws.Range("H132").Value = 1379.5714
$ws.Range("I132").Value = 1379.5714
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4138.7142
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1608.7142
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5412.5557
$ws.Range("J7").Value = 5464.25
$ws.Range("L7").Value = 5464.25
$ws.Range("N7").Value = -5688.25

$ws.Range("H40").Value = 15560606
$ws.Range("I40").Value = 6214.7144
$ws.Range("J40").Value = 51854188
$ws.Range("K40").Value = 6214.7144
$ws.Range("L40").Value = 51854188
$ws.Range("M40").Value = -6078.7144
$ws.Range("N40").Value = -51854460

$ws.Range("H68").Value = 6875.25
$ws.Range("I68").Value = 5833.6665
$ws.Range("K68").Value = 5833.6665
$ws.Range("M68").Value = -5084.6665

$ws.Range("H71").Value = 6875.25
$ws.Range("I71").Value = 5833.6665
$ws.Range("K71").Value = 29168.3325
$ws.Range("M71").Value = -25424.3325

$ws.Range("H122").Value = 120004060
$ws.Range("I122").Value = 142861150
$ws.Range("K122").Value = 428583450
$ws.Range("M122").Value = -428581000

$ws.Range("H126").Value = 5412.5557
$ws.Range("J126").Value = 5464.25
$ws.Range("L126").Value = 16392.75
$ws.Range("N126").Value = -21332.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 13249.75
$ws.Range("J69").Value = 13249.75
$ws.Range("L69").Value = 13249.75
$ws.Range("N69").Value = -14747.75

$ws.Range("H72").Value = 13249.75
$ws.Range("J72").Value = 13249.75
$ws.Range("L72").Value = 39749.25
$ws.Range("N72").Value = -47237.25

$ws.Range("H123").Value = 74800
$ws.Range("J123").Value = 74800
$ws.Range("L123").Value = 74800
$ws.Range("N123").Value = -84600

$ws.Range("H132").Value = 3240
$ws.Range("I132").Value = 3200
$ws.Range("K132").Value = 9600
$ws.Range("M132").Value = -7070

$ws.Range("H133").Value = 76497.664
$ws.Range("J133").Value = 76996.5
$ws.Range("L133").Value = 76996.5
$ws.Range("N133").Value = -87116.5
